# Applies the HEXACO fb report table edits:
#  1. "平均との比較" -> "値とその傾向" (and add w:hint="eastAsia" to that run's rFonts)
#  2-5. Merge the adjacent "<letter>" + "]," runs (e.g. "C" + "],") into a
#       single run with text "<letter>]," for each of C / E / A / N.

$d = $word.ActiveDocument

# --- 1. 平均との比較 -> 値とその傾向 -------------------------------------
# A plain Find/Replace changes the text but Word leaves the rFonts alone,
# so we replace the whole paragraph's XML (via the found Range) to match
# both the new text and the new w:hint="eastAsia" attribute exactly.
$r1 = $d.Content
$r1.Find.Execute("平均との比較", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r1.Find.Found) {
    $newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2758ACB9" w14:textId="77777777" w:rsidR="006863F8" w:rsidRPr="009D7D62" w:rsidRDefault="006863F8" w:rsidP="009D7D62"><w:pPr><w:widowControl w:val="0"/><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:jc w:val="center"/><w:rPr><w:rFonts w:asciiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r w:rsidRPr="009D7D62"><w:rPr><w:rFonts w:asciiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia" w:cs="Arial Unicode MS" w:hint="eastAsia"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>値とその傾向</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r1.InsertXML($newParaXml)
}

# --- 2-5. Merge "<letter>" + "]," runs --------------------------------
# These single-letter / "]," runs share identical rPr, so re-running
# Find/Replace over the already-correct text ("C]," -> "C],", etc.)
# causes Word to collapse them into one run instead of leaving a split,
# matching the diff exactly (keeps the first run's rPr, drops the second).
foreach ($letter in @("C", "E", "A", "N")) {
    $needle = "$letter],"
    $rr = $d.Content
    $rr.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
}
